$wb = $excel.ActiveWorkbook

# --- Rename sheets (update task-order identifiers embedded in sheet names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-165047781722466"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778189396904"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778189456599"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778190046597"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778190676632"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650477817182663.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778172076595.csv"
$ws1.Range("B4").Value = "go_stims-16504778172086592.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778172236927.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650477817629662.csv"
$ws2.Range("B3").Value = "OB-16504778177086885.csv"
$ws2.Range("B4").Value = "TB-165047781884566.csv"
$ws2.Range("B5").Value = "ZB-match_1-16504778174436626.csv"
$ws2.Range("B6").Value = "TB-16504778187426584.csv"
$ws2.Range("B7").Value = "ZB-match_0-16504778172876892.csv"
$ws2.Range("B8").Value = "ZB-match_9-16504778175556595.csv"
$ws2.Range("B9").Value = "TB-16504778189166949.csv"
$ws2.Range("B10").Value = "OB-16504778175816617.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650477818971691.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778189466588.csv"
$ws4.Range("B4").Value = "MM_stims-165047781898769.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477818972662.csv"
$ws4.Range("B6").Value = "MM_stims-165047781900369.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477818988662.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650477819035691.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778190196626.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778190516596.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778190076606.csv"
